$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 556.7143
$ws.Range("I6").Value = 249.4
$ws.Range("K6").Value = 748.2
$ws.Range("M6").Value = -636.2

$ws.Range("H32").Value = 912.2143
$ws.Range("I32").Value = 546
$ws.Range("J32").Value = 1115.6666
$ws.Range("K32").Value = 546
$ws.Range("L32").Value = 1115.6666
$ws.Range("M32").Value = -220
$ws.Range("N32").Value = -1767.6666

$ws.Range("H33").Value = 3636506.5
$ws.Range("I33").Value = 143.75
$ws.Range("J33").Value = 18181958
$ws.Range("K33").Value = 143.75
$ws.Range("L33").Value = 18181958
$ws.Range("M33").Value = 85.25
$ws.Range("N33").Value = -18182416

$ws.Range("H80").Value = 565.2083
$ws.Range("I80").Value = 745.5
$ws.Range("J80").Value = 529.15
$ws.Range("K80").Value = 2236.5
$ws.Range("L80").Value = 1587.45
$ws.Range("M80").Value = -1238.5
$ws.Range("N80").Value = -3583.45

$ws.Range("H83").Value = 565.2083
$ws.Range("I83").Value = 745.5
$ws.Range("J83").Value = 529.15
$ws.Range("K83").Value = 6709.5
$ws.Range("L83").Value = 4762.349999999999
$ws.Range("M83").Value = -1717.5
$ws.Range("N83").Value = -14746.35

$ws.Range("H113").Value = 6407.5654
$ws.Range("I113").Value = 4760.9
$ws.Range("J113").Value = 7674.231
$ws.Range("K113").Value = 4760.9
$ws.Range("L113").Value = 7674.231
$ws.Range("M113").Value = -1506.9
$ws.Range("N113").Value = -14182.231

$ws.Range("H126").Value = 29197.777
$ws.Range("J126").Value = 29197.777
$ws.Range("L126").Value = 29197.777
$ws.Range("N126").Value = -39077.777

$ws.Range("H138").Value = 2600.407
$ws.Range("J138").Value = 3840.8542
$ws.Range("L138").Value = 11522.5626
$ws.Range("N138").Value = -21802.5626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4000.5571
$ws.Range("I32").Value = 3354.7114
$ws.Range("J32").Value = 5866.3335
$ws.Range("K32").Value = 3354.7114
$ws.Range("L32").Value = 5866.3335
$ws.Range("M32").Value = -3067.7114
$ws.Range("N32").Value = -6440.3335

$ws.Range("H61").Value = 5915.231
$ws.Range("I61").Value = 2375.3333
$ws.Range("J61").Value = 20782.8
$ws.Range("K61").Value = 2375.3333
$ws.Range("L61").Value = 20782.8
$ws.Range("M61").Value = -2163.3333
$ws.Range("N61").Value = -21206.8

$ws.Range("H136").Value = 5915.231
$ws.Range("I136").Value = 2375.3333
$ws.Range("J136").Value = 20782.8
$ws.Range("K136").Value = 7125.999899999999
$ws.Range("L136").Value = 62348.39999999999
$ws.Range("M136").Value = -4575.999899999999
$ws.Range("N136").Value = -67448.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 142859310
$ws.Range("I20").Value = 500001250
$ws.Range("J20").Value = 2538
$ws.Range("K20").Value = 500001250
$ws.Range("L20").Value = 2538
$ws.Range("M20").Value = -500001003
$ws.Range("N20").Value = -3032

$ws.Range("H22").Value = 330.25
$ws.Range("I22").Value = 330.25
$ws.Range("K22").Value = 330.25
$ws.Range("M22").Value = -157.25

$ws.Range("H107").Value = 1800
$ws.Range("I107").Value = 1800
$ws.Range("K107").Value = 1800
$ws.Range("M107").Value = 120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2589.0137
$ws.Range("I31").Value = 2135.389
$ws.Range("J31").Value = 2737.4727
$ws.Range("K31").Value = 2135.389
$ws.Range("L31").Value = 2737.4727
$ws.Range("M31").Value = -1840.389
$ws.Range("N31").Value = -3327.4727

$ws.Range("H34").Value = 2589.0137
$ws.Range("I34").Value = 2135.389
$ws.Range("J34").Value = 2737.4727
$ws.Range("K34").Value = 2135.389
$ws.Range("L34").Value = 2737.4727
$ws.Range("M34").Value = -1933.389
$ws.Range("N34").Value = -3141.4727

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0

$ws.Range("H99").Value = 11380427
$ws.Range("I99").Value = 15483.167
$ws.Range("J99").Value = 25018360
$ws.Range("K99").Value = 15483.167
$ws.Range("L99").Value = 25018360
$ws.Range("M99").Value = -13985.167
$ws.Range("N99").Value = -25021356

$ws.Range("H126").Value = 11380427
$ws.Range("I126").Value = 15483.167
$ws.Range("J126").Value = 25018360
$ws.Range("K126").Value = 46449.501
$ws.Range("L126").Value = 75055080
$ws.Range("M126").Value = -43979.501
$ws.Range("N126").Value = -75060020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 149680.64
$ws.Range("I5").Value = 14790.857
$ws.Range("J5").Value = 177452.06
$ws.Range("K5").Value = 44372.571
$ws.Range("L5").Value = 532356.1799999999
$ws.Range("M5").Value = -44260.571
$ws.Range("N5").Value = -532580.1799999999

$ws.Range("H7").Value = 1029
$ws.Range("I7").Value = 257.5
$ws.Range("J7").Value = 1800.5
$ws.Range("K7").Value = 772.5
$ws.Range("L7").Value = 5401.5
$ws.Range("M7").Value = -660.5
$ws.Range("N7").Value = -5625.5

$ws.Range("H68").Value = 3595.8462
$ws.Range("J68").Value = 2376.6155
$ws.Range("L68").Value = 7129.8465
$ws.Range("N68").Value = -8751.8465

$ws.Range("H71").Value = 3595.8462
$ws.Range("J71").Value = 2376.6155
$ws.Range("L71").Value = 21389.5395
$ws.Range("N71").Value = -29501.5395

$ws.Range("H92").Value = 758
$ws.Range("I92").Value = 760.6
$ws.Range("J92").Value = 751.5
$ws.Range("K92").Value = 2281.8
$ws.Range("L92").Value = 2254.5
$ws.Range("M92").Value = -1033.8
$ws.Range("N92").Value = -4750.5

$ws.Range("H135").Value = 149680.64
$ws.Range("I135").Value = 14790.857
$ws.Range("J135").Value = 177452.06
$ws.Range("K135").Value = 133117.713
$ws.Range("L135").Value = 1597068.54
$ws.Range("M135").Value = -130582.713
$ws.Range("N135").Value = -1602138.54

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1683.25
$ws.Range("I102").Value = 1344.3334
$ws.Range("K102").Value = 1344.3334
$ws.Range("M102").Value = 277.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2851165.2
$ws.Range("I22").Value = 15874145
$ws.Range("J22").Value = 2388.5625
$ws.Range("K22").Value = 15874145
$ws.Range("L22").Value = 2388.5625
$ws.Range("M22").Value = -15873850
$ws.Range("N22").Value = -2978.5625

$ws.Range("H27").Value = 2851165.2
$ws.Range("I27").Value = 15874145
$ws.Range("J27").Value = 2388.5625
$ws.Range("K27").Value = 15874145
$ws.Range("L27").Value = 2388.5625
$ws.Range("M27").Value = -15874038
$ws.Range("N27").Value = -2602.5625

$ws.Range("H55").Value = 20833798
$ws.Range("I55").Value = 444.33334
$ws.Range("J55").Value = 33333810
$ws.Range("K55").Value = 444.33334
$ws.Range("L55").Value = 33333810
$ws.Range("M55").Value = -271.33334
$ws.Range("N55").Value = -33334156

$ws.Range("H61").Value = 2899.5715
$ws.Range("I61").Value = 2882.8333
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2882.8333
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2680.8333
$ws.Range("N61").Value = -3404

$ws.Range("H80").Value = 32000
$ws.Range("J80").Value = 32000
$ws.Range("L80").Value = 32000
$ws.Range("N80").Value = -34246

$ws.Range("H83").Value = 32000
$ws.Range("J83").Value = 32000
$ws.Range("L83").Value = 96000
$ws.Range("N83").Value = -107232

$ws.Range("H113").Value = 2899.5715
$ws.Range("I113").Value = 2882.8333
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2882.8333
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -712.8332999999998
$ws.Range("N113").Value = -7340

$ws.Range("H122").Value = 7407597.5
$ws.Range("I122").Value = 14289714
$ws.Range("K122").Value = 42869142
$ws.Range("M122").Value = -42866692

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7062.3076
$ws.Range("J54").Value = 7062.3076
$ws.Range("L54").Value = 7062.3076
$ws.Range("N54").Value = -8102.3076

$ws.Range("H126").Value = 1375.0625
$ws.Range("I126").Value = 936.44446
$ws.Range("K126").Value = 2809.33338
$ws.Range("M126").Value = -339.33338
